# Added stock update through socket io problem fixed
#
# Inserts a new "stock_quantity" column (AW) into the products sheet,
# shifting the former AW (ID) -> AX and AX (date_created) -> AY for the
# rows that already had those trailing ID/date_created columns (18-27).
# Also updates a handful of other cell values (size/volume labels,
# category/type corrections, stock counts, and last-modified timestamps).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---------------------------------------------------------------------
# 1) Header row: insert new "stock_quantity" header at AW1, push the old
#    AW1 ("ID") to AX1 and old AX1 ("date_created") to AY1.
# ---------------------------------------------------------------------
$ws.Range("AY1").Value = "date_created"
$ws.Range("AX1").Value = "ID"
$ws.Range("AW1").Value = "stock_quantity"

# ---------------------------------------------------------------------
# 2) Row 2 (Dom Perignon Vintage Champagne6) updates
# ---------------------------------------------------------------------
$ws.Range("AB2").Value = "7ml"
$ws.Range("AC2").Value = "7ML"
$ws.Range("AG2").Value = 1
$ws.Range("AU2").Value = "2025-03-31T17:11:56.419Z"
$ws.Range("AW2").Value = 1

# ---------------------------------------------------------------------
# 3) Row 3 updates (scotch -> vodka recategorisation, size fix)
# ---------------------------------------------------------------------
$ws.Range("G3").Value = "vodka"
$ws.Range("H3").Value = "vodka"
$ws.Range("AB3").Value = "75ml"
$ws.Range("AC3").Value = "75ML"
$ws.Range("AI3").Value = "vodka"
$ws.Range("AU3").Value = "2025-03-31T17:01:46.261Z"
$ws.Range("AW3").Value = 12

# ---------------------------------------------------------------------
# 4) Rows 4-11: just add the new stock_quantity value
# ---------------------------------------------------------------------
$ws.Range("AW4").Value = 32
$ws.Range("AW5").Value = 34
$ws.Range("AW6").Value = 45
$ws.Range("AW7").Value = 45
$ws.Range("AW8").Value = 65
$ws.Range("AW9").Value = 65
$ws.Range("AW10").Value = 66
$ws.Range("AW11").Value = 66

# ---------------------------------------------------------------------
# 5) Row 12 updates
# ---------------------------------------------------------------------
$ws.Range("AB12").Value = "63ml"
$ws.Range("AC12").Value = "63ML"
$ws.Range("AG12").Value = 9
$ws.Range("AU12").Value = "2025-03-31T17:08:09.075Z"
$ws.Range("AW12").Value = 9

# ---------------------------------------------------------------------
# 6) Rows 13-17: just add the new stock_quantity value
# ---------------------------------------------------------------------
$ws.Range("AW13").Value = 6
$ws.Range("AW14").Value = 5
$ws.Range("AW15").Value = 4
$ws.Range("AW16").Value = 45
$ws.Range("AW17").Value = 45

# ---------------------------------------------------------------------
# 7) Rows 18-27: these rows already had trailing "ID"/"date_created"
#    values stored (as text) in AW/AX. Shift them right to AX/AY and
#    insert the new numeric stock_quantity into AW. Row 25 never had an
#    ID/date_created pair, so it only gets the new stock_quantity value.
#    The "'"-prefix forces these numeric-looking IDs to stay text cells,
#    matching the source file's t="str" cells.
# ---------------------------------------------------------------------
$ws.Range("AW18").Value = 32
$ws.Range("AX18").Value = "'600"
$ws.Range("AY18").Value = "2025-03-28T12:17:35.575Z"

$ws.Range("AW19").Value = 2
$ws.Range("AX19").Value = "'601"
$ws.Range("AY19").Value = "2025-03-28T12:25:00.293Z"

$ws.Range("AW20").Value = 3
$ws.Range("AX20").Value = "'602"
$ws.Range("AY20").Value = "2025-03-28T12:31:41.045Z"

$ws.Range("AW21").Value = 4
$ws.Range("AX21").Value = "'603"
$ws.Range("AY21").Value = "2025-03-28T12:46:29.945Z"

$ws.Range("AW22").Value = 5
$ws.Range("AX22").Value = "'604"
$ws.Range("AY22").Value = "2025-03-28T13:16:01.109Z"

$ws.Range("AW23").Value = 6
$ws.Range("AX23").Value = "'605"
$ws.Range("AY23").Value = "2025-03-28T14:45:57.115Z"

$ws.Range("AW24").Value = 7
$ws.Range("AX24").Value = "'606"
$ws.Range("AY24").Value = "2025-03-28T15:15:57.426Z"

$ws.Range("AW25").Value = 6

$ws.Range("AW26").Value = 7
$ws.Range("AX26").Value = "'608"
$ws.Range("AY26").Value = "2025-03-29T05:34:46.251Z"

$ws.Range("AW27").Value = 9
$ws.Range("AX27").Value = "'609"
$ws.Range("AY27").Value = "2025-03-29T09:07:41.265Z"

$wb.Save()
